{"js": "// Wrap each of the four bullet-point paragraphs (report summary items)\n// with a bookmark, matching the target diff:\n//   \"Quantidade de projetos por status\"          -> __DdeLink__114_1841298793\n//   \"Total or\u00e7ado por status\"                    -> __DdeLink__116_1841298793\n//   \"M\u00e9dia de dura\u00e7\u00e3o dos projetos encerrados\"   -> __DdeLink__118_1841298793\n//   \"Total de membros \u00fanicos alocados\"           -> __DdeLink__120_1841298793\n\nconst targets = [\n  { text: \"Quantidade de projetos por status\", name: \"__DdeLink__114_1841298793\" },\n  { text: \"Total or\u00e7ado por status\", name: \"__DdeLink__116_1841298793\" },\n  { text: \"M\u00e9dia de dura\u00e7\u00e3o dos projetos encerrados\", name: \"__DdeLink__118_1841298793\" },\n  { text: \"Total de membros \u00fanicos alocados\", name: \"__DdeLink__120_1841298793\" },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const target of targets) {\n  const paragraph = paragraphs.items.find((p) => p.text.includes(target.text));\n  if (!paragraph) {\n    throw new Error(`Paragraph containing \"${target.text}\" not found`);\n  }\n  // \"Content\" keeps the bookmark scoped to the paragraph's runs (not the\n  // trailing paragraph mark), matching bookmarkStart/bookmarkEnd placed\n  // right after <w:pPr> and right before </w:p>.\n  paragraph.getRange(\"Content\").insertBookmark(target.name);\n}\n\nawait context.sync();\n", "ps1": "# Wrap each of the four bullet-point paragraphs (report summary items)\n# with a bookmark, matching the target diff:\n#   \"Quantidade de projetos por status\"          -> __DdeLink__114_1841298793\n#   \"Total or\u00e7ado por status\"                    -> __DdeLink__116_1841298793\n#   \"M\u00e9dia de dura\u00e7\u00e3o dos projetos encerrados\"   -> __DdeLink__118_1841298793\n#   \"Total de membros \u00fanicos alocados\"           -> __DdeLink__120_1841298793\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    @{ Text = \"Quantidade de projetos por status\"; Name = \"__DdeLink__114_1841298793\" },\n    @{ Text = \"Total or\u00e7ado por status\"; Name = \"__DdeLink__116_1841298793\" },\n    @{ Text = \"M\u00e9dia de dura\u00e7\u00e3o dos projetos encerrados\"; Name = \"__DdeLink__118_1841298793\" },\n    @{ Text = \"Total de membros \u00fanicos alocados\"; Name = \"__DdeLink__120_1841298793\" }\n)\n\nforeach ($target in $targets) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -like \"*$($target.Text)*\") {\n            $r = $p.Range\n            # Exclude the trailing paragraph mark so the bookmark wraps only\n            # the paragraph's own content (matches bookmarkStart right after\n            # <w:pPr> and bookmarkEnd right before </w:p>).\n            $r.End = $r.End - 1\n            $d.Bookmarks.Add($target.Name, $r)\n            break\n        }\n    }\n}\n"}
